$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 15560.2
$ws.Range("J87").Value = 15560.2
$ws.Range("L87").Value = 15560.2
$ws.Range("N87").Value = -18056.2
$ws.Range("H90").Value = 15560.2
$ws.Range("J90").Value = 15560.2
$ws.Range("L90").Value = 46680.60000000001
$ws.Range("N90").Value = -59160.60000000001
$ws.Range("H129").Value = 833.2838
$ws.Range("J129").Value = 890.1875
$ws.Range("L129").Value = 2670.5625
$ws.Range("N129").Value = -12670.5625
$ws.Range("H137").Value = 889.871
$ws.Range("I137").Value = 830.43475
$ws.Range("J137").Value = 1060.75
$ws.Range("K137").Value = 2491.30425
$ws.Range("L137").Value = 3182.25
$ws.Range("M137").Value = 58.69574999999986
$ws.Range("N137").Value = -8282.25
$ws.Range("H138").Value = 3492.65
$ws.Range("J138").Value = 4788.393
$ws.Range("L138").Value = 14365.179
$ws.Range("N138").Value = -24645.179

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1727.4736
$ws.Range("I74").Value = 942.2727
$ws.Range("K74").Value = 942.2727
$ws.Range("M74").Value = -68.27269999999999
$ws.Range("H77").Value = 1727.4736
$ws.Range("I77").Value = 942.2727
$ws.Range("K77").Value = 4711.363499999999
$ws.Range("M77").Value = -343.3634999999995
$ws.Range("H122").Value = 2532.5
$ws.Range("I122").Value = 1711.25
$ws.Range("J122").Value = 4175
$ws.Range("K122").Value = 5133.75
$ws.Range("L122").Value = 12525
$ws.Range("M122").Value = -2683.75
$ws.Range("N122").Value = -17425
$ws.Range("H132").Value = 1763.2778
$ws.Range("I132").Value = 1407.8918
$ws.Range("J132").Value = 2536.7646
$ws.Range("K132").Value = 4223.6754
$ws.Range("L132").Value = 7610.293799999999
$ws.Range("M132").Value = -1693.6754
$ws.Range("N132").Value = -12670.2938

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1288.6177
$ws.Range("I134").Value = 1050.7916
$ws.Range("J134").Value = 1859.4
$ws.Range("K134").Value = 3152.3748
$ws.Range("L134").Value = 5578.200000000001
$ws.Range("M134").Value = -617.3748000000001
$ws.Range("N134").Value = -10648.2
$ws.Range("H140").Value = 54844.285
$ws.Range("J140").Value = 54844.285
$ws.Range("L140").Value = 54844.285
$ws.Range("N140").Value = -65204.285

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2516.24
$ws.Range("I31").Value = 1612.6545
$ws.Range("J31").Value = 3620.6223
$ws.Range("K31").Value = 1612.6545
$ws.Range("L31").Value = 3620.6223
$ws.Range("M31").Value = -1317.6545
$ws.Range("N31").Value = -4210.6223
$ws.Range("H34").Value = 2516.24
$ws.Range("I34").Value = 1612.6545
$ws.Range("J34").Value = 3620.6223
$ws.Range("K34").Value = 1612.6545
$ws.Range("L34").Value = 3620.6223
$ws.Range("M34").Value = -1410.6545
$ws.Range("N34").Value = -4024.6223
$ws.Range("H58").Value = 3122.7307
$ws.Range("I58").Value = 2773.5557
$ws.Range("J58").Value = 3307.5881
$ws.Range("K58").Value = 2773.5557
$ws.Range("L58").Value = 3307.5881
$ws.Range("M58").Value = -2570.5557
$ws.Range("N58").Value = -3713.5881
$ws.Range("H99").Value = 1785786.4
$ws.Range("I99").Value = 2135410.2
$ws.Range("J99").Value = 37666.668
$ws.Range("K99").Value = 2135410.2
$ws.Range("L99").Value = 37666.668
$ws.Range("M99").Value = -2133912.2
$ws.Range("N99").Value = -40662.668
$ws.Range("H126").Value = 1785786.4
$ws.Range("I126").Value = 2135410.2
$ws.Range("J126").Value = 37666.668
$ws.Range("K126").Value = 6406230.600000001
$ws.Range("L126").Value = 113000.004
$ws.Range("M126").Value = -6403760.600000001
$ws.Range("N126").Value = -117940.004
$ws.Range("H136").Value = 3122.7307
$ws.Range("I136").Value = 2773.5557
$ws.Range("J136").Value = 3307.5881
$ws.Range("K136").Value = 8320.667099999999
$ws.Range("L136").Value = 9922.764299999999
$ws.Range("M136").Value = -5770.667099999999
$ws.Range("N136").Value = -15022.7643
$ws.Range("H140").Value = 82951.164
$ws.Range("J140").Value = 82951.164
$ws.Range("L140").Value = 82951.164
$ws.Range("N140").Value = -93311.164

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 924.62
$ws.Range("I131").Value = 665
$ws.Range("J131").Value = 929.9184
$ws.Range("K131").Value = 1995
$ws.Range("L131").Value = 2789.7552
$ws.Range("M131").Value = 3045
$ws.Range("N131").Value = -12869.7552

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4527461.5
$ws.Range("I102").Value = 5130456.5
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 5130456.5
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -5128834.5
$ws.Range("N102").Value = -8244
$ws.Range("H113").Value = 1492.4736
$ws.Range("I113").Value = 1464.25
$ws.Range("J113").Value = 1540.8572
$ws.Range("K113").Value = 1464.25
$ws.Range("L113").Value = 1540.8572
$ws.Range("M113").Value = 705.75
$ws.Range("N113").Value = -5880.8572
$ws.Range("H132").Value = 2394.9656
$ws.Range("I132").Value = 2022.1
$ws.Range("J132").Value = 3223.5557
$ws.Range("K132").Value = 6066.299999999999
$ws.Range("L132").Value = 9670.667099999999
$ws.Range("M132").Value = -3536.299999999999
$ws.Range("N132").Value = -14730.6671
$ws.Range("H138").Value = 48250.715
$ws.Range("J138").Value = 48250.715
$ws.Range("L138").Value = 48250.715
$ws.Range("N138").Value = -58530.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 84899.836
$ws.Range("I40").Value = 84899.836
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 84899.836
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -84763.836
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1503.95
$ws.Range("I46").Value = 1205.6428
$ws.Range("J46").Value = 2200
$ws.Range("K46").Value = 1205.6428
$ws.Range("L46").Value = 2200
$ws.Range("M46").Value = -1017.6428
$ws.Range("N46").Value = -2576
$ws.Range("H55").Value = 215.76471
$ws.Range("I55").Value = 221.125
$ws.Range("J55").Value = 211
$ws.Range("K55").Value = 221.125
$ws.Range("L55").Value = 211
$ws.Range("M55").Value = -48.125
$ws.Range("N55").Value = -557
$ws.Range("H100").Value = 3117.875
$ws.Range("I100").Value = 1489.8334
$ws.Range("J100").Value = 8002
$ws.Range("K100").Value = 1489.8334
$ws.Range("L100").Value = 8002
$ws.Range("M100").Value = -948.8334
$ws.Range("N100").Value = -9084
$ws.Range("H127").Value = 39791
$ws.Range("J127").Value = 39791
$ws.Range("L127").Value = 39791
$ws.Range("N127").Value = -49711
$ws.Range("H136").Value = 25644856
$ws.Range("I136").Value = 4128.8
$ws.Range("J136").Value = 111113944
$ws.Range("K136").Value = 12386.4
$ws.Range("L136").Value = 333341832
$ws.Range("M136").Value = -9836.400000000001
$ws.Range("N136").Value = -333346932

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1539.5074
$ws.Range("I132").Value = 977.2439000000001
$ws.Range("J132").Value = 2426.1538
$ws.Range("K132").Value = 2931.7317
$ws.Range("L132").Value = 7278.4614
$ws.Range("M132").Value = -401.7317000000003
$ws.Range("N132").Value = -12338.4614
$ws.Range("H137").Value = 55857.168
$ws.Range("J137").Value = 55857.168
$ws.Range("L137").Value = 55857.168
$ws.Range("N137").Value = -66057.16800000001
